# Update the "empadronador" / "total_registros" table with refreshed
# monitoring data. The rows are rewritten (both names and counts) and
# re-sorted by total_registros descending.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("RUIZ CHIROQUE CLAUDIA JUDITH", 130),
    @("GONZALES FIESTAS MARIA MARIBEL", 124),
    @("BANCAYAN FIESTA DILVER HUMBERTO", 122),
    @("TEMOCHE ECHE URSULA YESSENIA", 121),
    @("LLENQUE ANTON HELEN JOHANA", 120),
    @("FABIANA REBECA ARRUNATEGUI SILUPU", 120),
    @("ANTON INGA FATIMA DEL ROSARIO", 119),
    @("BAUTISTA CHAVESTA ERICKA MEDALIT", 105),
    @("FLORES SILUPU MARY CARMEN", 101),
    @("VELASCO PEÑA KAREN ARELLYS", 101),
    @("PINTADO CHASQUERO ESTEFANY", 100),
    @("HERNANDEZ CARNERO ARTURO SEBASTIAN", 95),
    @("MONDRAGON NONAJULCA MARISOL", 93),
    @("ORDINOLA JIBAJA JOSE ALBERTO", 91),
    @("CASTRO ESTRADA CINTHIA PATRICIA", 71),
    @("MORENO YANAYACO NAYLA GUADALUPE", 66),
    @("PINTADO BENITES CRISTOBAL RODRIGO", 1)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $row = $row + 1
}
